$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was updated
# from serial 45186 (2023-09-17) to serial 45188 (2023-09-19) for
# every data row (rows 2 through 244).
$newDate = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).Date.AddDays(45188)

for ($r = 2; $r -le 244; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}
